$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 10 rows before the existing data (rows 2-11) to make room for a new "Brasil" data block.
# This shifts the existing "Nordeste" rows down from (2-11) to (12-21) and "Sergipe" rows down from (12-21) to (22-31).
$ws.Rows("2:11").Insert()

# The inserted rows pick up formatting from the row above (the header); clear it so the new rows
# look like ordinary (unstyled) data rows, consistent with the rest of the table.
$ws.Rows("2:11").ClearFormats()

# Fill in the newly inserted rows with the "Brasil" data series.
$ws.Cells.Item(2, 1).Value = "Brasil"
$ws.Cells.Item(2, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(2, 3).Value = "31/12/2013"
$ws.Cells.Item(2, 4).Value = 100

$ws.Cells.Item(3, 1).Value = "Brasil"
$ws.Cells.Item(3, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(3, 3).Value = "31/12/2014"
$ws.Cells.Item(3, 4).Value = 101.2729253253112

$ws.Cells.Item(4, 1).Value = "Brasil"
$ws.Cells.Item(4, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(4, 3).Value = "31/12/2015"
$ws.Cells.Item(4, 4).Value = 99.86379134956168

$ws.Cells.Item(5, 1).Value = "Brasil"
$ws.Cells.Item(5, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(5, 3).Value = "31/12/2016"
$ws.Cells.Item(5, 4).Value = 94.09943317286582

$ws.Cells.Item(6, 1).Value = "Brasil"
$ws.Cells.Item(6, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(6, 3).Value = "31/12/2017"
$ws.Cells.Item(6, 4).Value = 94.55172957222145

$ws.Cells.Item(7, 1).Value = "Brasil"
$ws.Cells.Item(7, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(7, 3).Value = "31/12/2018"
$ws.Cells.Item(7, 4).Value = 95.26579737496398

$ws.Cells.Item(8, 1).Value = "Brasil"
$ws.Cells.Item(8, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(8, 3).Value = "31/12/2019"
$ws.Cells.Item(8, 4).Value = 103.3813605432476

$ws.Cells.Item(9, 1).Value = "Brasil"
$ws.Cells.Item(9, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(9, 3).Value = "31/12/2020"
$ws.Cells.Item(9, 4).Value = 100.5542731061483

$ws.Cells.Item(10, 1).Value = "Brasil"
$ws.Cells.Item(10, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(10, 3).Value = "31/12/2021"
$ws.Cells.Item(10, 4).Value = 99.55144222900864

$ws.Cells.Item(11, 1).Value = "Brasil"
$ws.Cells.Item(11, 2).Value = "Índice do emprego formal: 2013=100"
$ws.Cells.Item(11, 3).Value = "31/12/2022"
$ws.Cells.Item(11, 4).Value = 107.8499571171155
